# Reorder the worksheet tabs: "review_info" should become the first sheet,
# "hotel_info" the second sheet.
$wb = $excel.ActiveWorkbook
$reviewInfo = $wb.Worksheets.Item("review_info")
$reviewInfo.Move($wb.Worksheets.Item(1))

# Re-fetch hotel_info by name now that the sheets have been reordered.
$hotelInfo = $wb.Worksheets.Item("hotel_info")

# Insert a new "State" column into hotel_info between "Hotel_Name" and "City",
# populated with "Louisiana" for the existing data row.
$hotelInfo.Columns.Item(3).Insert()
$hotelInfo.Cells.Item(1, 3).Value = "State"
$hotelInfo.Cells.Item(2, 3).Value = "Louisiana"
